# Update "paises.xlsx" — refresh country case counts and reorder three
# country labels (Belice / Nueva Caledonia / Santa Lucia), matching the
# "Update countries & provincias Spain" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp header (A1) ------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 11:05"

# --- Refresh numeric data for a few countries --------------------------
# Row 19 = Belgica
$ws.Range("B19").Value = 55983
$ws.Range("C19").Value = 192
$ws.Range("D19").Value = 14847
$ws.Range("E19").Value = 31986
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 9150

# Row 31 = Banglades
$ws.Range("B31").Value = 26738
$ws.Range("C31").Value = 1617
$ws.Range("D31").Value = 5207
$ws.Range("E31").Value = 21145
$ws.Range("G31").Value = 16
$ws.Range("H31").Value = 386

# Row 59 = Malasia
$ws.Range("B59").Value = 7009
$ws.Range("C59").Value = 31
$ws.Range("D59").Value = 5706
$ws.Range("E59").Value = 1189

# Row 200 = Namibia
$ws.Range("D200").Value = 14
$ws.Range("E200").Value = 2

# --- Reorder Belice / Nueva Caledonia / Santa Lucia ---------------------
# Rows 195-197 keep their position but the labels rotate: the row that
# used to read "Belice" now reads "Nueva Caledonia", the one that used to
# read "Nueva Caledonia" now reads "Santa Lucia", and the one that used to
# read "Santa Lucia" now reads "Belice" — with their data following the
# label.
$ws.Range("A195").Value = "Nueva Caledonia"
$ws.Range("D195").Value = 18
$ws.Range("H195").Value = 0

$ws.Range("A196").Value = "Santa Lucia"

$ws.Range("A197").Value = "Belice"
$ws.Range("D197").Value = 16
$ws.Range("H197").Value = 2
